# "images on the end of month"
# Append the September 2025 trading-day PE/PB observations to the myPEPB
# sheet, continuing the running row-count/cumulative-average formulas that
# were already in place through row 1072 (2025/8/29).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)   # myPEPB

$startRow = 1073
$endRow = 1094

# Give the new rows the same look as the existing data block: column B in
# the yyyy/mm/dd date format, column C with 3 decimal places, both right
# aligned (matching styles already used through row 1072); A/D stay
# general, same as above.
$ws.Range("B$startRow`:B$endRow").NumberFormat = "yyyy/mm/dd"
$ws.Range("B$startRow`:B$endRow").HorizontalAlignment = -4152
$ws.Range("C$startRow`:C$endRow").NumberFormat = "0.000_ "
$ws.Range("C$startRow`:C$endRow").HorizontalAlignment = -4152

# Column B alternates between text "yyyy/m/d`n" values (pasted from an
# external data source, like the bulk of the sheet already does) and, for
# the last two days of the month, plain date serials -- exactly mirroring
# the pattern already present at the end of every previous month's block
# (see rows 1071/1072, 1050/1051, ...).
$dates = @(
    "2025/9/1`n",
    "2025/9/2`n",
    "2025/9/3`n",
    "2025/9/4`n",
    "2025/9/5`n",
    "2025/9/8`n",
    "2025/9/9`n",
    "2025/9/10`n",
    "2025/9/11`n",
    "2025/9/12`n",
    "2025/9/15`n",
    "2025/9/16`n",
    "2025/9/17`n",
    "2025/9/18`n",
    "2025/9/19`n",
    "2025/9/22`n",
    "2025/9/23`n",
    "2025/9/24`n",
    "2025/9/25`n",
    "2025/9/26`n",
    45929,
    45930
)

$values = @(
    29.75,
    29.120000839999999,
    29.170000080000001,
    28.18000031,
    28.18000031,
    29.68000031,
    29.090000150000002,
    29.219999309999999,
    30.219999309999999,
    30,
    30.559999470000001,
    30.760000229999999,
    31.459999079999999,
    31.329999919999999,
    31.159999849999998,
    31.489999770000001,
    31.520000459999999,
    32.240001679999999,
    32.72000122,
    32.009998320000001,
    32.75,
    32.849998470000003
)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $r = $startRow + $i
    $prev = $r - 1

    $ws.Cells.Item($r, 1).Formula = "=A$prev+1"
    $ws.Cells.Item($r, 2).Value = $dates[$i]
    $ws.Cells.Item($r, 3).Value = $values[$i]
    $ws.Cells.Item($r, 4).Formula = "=SUM(C`$3:C$r)/A$r"
}

# The embedded newline in the text dates makes the headless host's
# auto-height kick in; put row heights back to the sheet default like the
# pre-existing rows above them.
$ws.Range("A$startRow`:D$endRow").EntireRow.AutoFit()
